$wb = $excel.ActiveWorkbook

# =========================================================================
# 1. Add the new "Subscription" worksheet as the last sheet in the workbook
# =========================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$subSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$subSheet.Name = "Subscription"

# =========================================================================
# 2. Update the "ContactUs" sheet
# =========================================================================
$contact = $wb.Worksheets.Item("ContactUs")

# Row 5 - update the expected message text
$contact.Range("E5").Value = "Please fill out this field."

# Row 7 - brand new test case row
$contact.Range("A7").Value = "John@123#"
$contact.Range("B7").Value = "validemail@example.com"
$contact.Range("C7").Value = "Invalid Name TC"
$contact.Range("D7").Value = "Testing invalid characters."
$contact.Range("E7").Value = "INVALID_NAME"

# Hyperlink for the email address cell (re-assert the display text afterwards so
# the cell keeps showing the plain address rather than the "mailto:" link text)
$contact.Hyperlinks.Add($contact.Range("B7"), "mailto:validemail@example.com", "", "", "mailto:validemail@example.com")
$contact.Range("B7").Value = "validemail@example.com"

# Formatting for row 7: wrap + vertical centered, taller row
$contact.Range("A7:E7").WrapText = $true
$contact.Range("A7:E7").VerticalAlignment = -4108
$contact.Rows.Item(7).RowHeight = 29

# Selection / active cell bookkeeping to mirror the saved workbook state
$contact.Range("B12").Select()

# =========================================================================
# 3. Populate the new "Subscription" sheet
# =========================================================================
$subSheet.Range("A1").Value = "Email"
$subSheet.Range("B1").Value = "ExpectedMessage"
$subSheet.Range("A2").Value = "valid@example.com"
$subSheet.Range("B2").Value = "You have been successfully subscribed!"
$subSheet.Range("B3").Value = "Please fill out this field"
$subSheet.Range("A4").Value = "invalidemail"
$subSheet.Range("B4").Value = "Please include an '@' in the email address"

# Hyperlink for the email address cell (re-assert the display text afterwards so
# the cell keeps showing the plain address rather than the "mailto:" link text)
$subSheet.Hyperlinks.Add($subSheet.Range("A2"), "mailto:valid@example.com", "", "", "mailto:valid@example.com")
$subSheet.Range("A2").Value = "valid@example.com"

# Header formatting: bold, centered, wrap text
$subSheet.Range("A1:B1").Font.Bold = $true
$subSheet.Range("A1:B1").WrapText = $true
$subSheet.Range("A1:B1").HorizontalAlignment = -4108
$subSheet.Range("A1:B1").VerticalAlignment = -4108

# Body formatting: wrap + vertical centered
$subSheet.Range("A2:B4").WrapText = $true
$subSheet.Range("A2:B4").VerticalAlignment = -4108
$subSheet.Range("A3").WrapText = $true
$subSheet.Range("A3").VerticalAlignment = -4108

$subSheet.Rows.Item(2).RowHeight = 29

# Column widths
$subSheet.Columns.Item(1).ColumnWidth = 12.893229166666666
$subSheet.Columns.Item(2).ColumnWidth = 41.529947916666664

$subSheet.Range("F2").Select()

# =========================================================================
# 4. Re-activate "ContactUs" so it remains the selected tab, as in the
#    original workbook.
# =========================================================================
$contact.Activate()
$contact.Range("B12").Select()
